$d = $word.ActiveDocument

# The cover-page version/date line reads "Version 11.02.03, 2015-05-27"
# and must become "Version 11.03.10, 2015-06-29" (commit: update docs to
# reflect the 11.03.10 release). The leading "Version 11.0" text is left
# untouched (it is outside the diff); only the trailing
# "2.03, 2015-05-27" -> "3.10, 2015-06-29" span is replaced.
$d.Content.Find.Execute("2.03, 2015-05-27", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.10, 2015-06-29", 2)
